# Generate Report for Handback
# The 45e61196-... file has now been handed back (in sync with en-US), so it
# moves to the top of each report table, and the 304848e7-... file (still
# "Ready for handoff") drops to the second row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 -> 45e61196 file, now handed back
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-25-11 22:25:38"

# Row 3 -> 304848e7 file, still ready for handoff
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-25-11 22:25:15"

# The hyperlinks on column A keep pointing at their original targets but the
# *displayed* file name swaps along with the row contents.
$ov.Range("A2").Hyperlinks.Delete()
$ov.Range("A3").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/304848e7-1a06-40e6-bfe6-b511954e83f4.md", "", "", "45e61196-dd2b-4863-a39e-77c67d4820a6.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3608f97c3b8b650493d45bd665d1432dfa3ab57c/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md", "", "", "304848e7-1a06-40e6-bfe6-b511954e83f4.md") | Out-Null

# ---------------------------------------------------------------------------
# Helper data shared by the two language detail sheets
# ---------------------------------------------------------------------------
function Fill-LangSheet($ws, $xlfExt) {
    $mdUrl45 = "https://github.com/OpenLocalizationTest/oltest/blob/3608f97c3b8b650493d45bd665d1432dfa3ab57c/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md"
    $mdUrl30 = "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/304848e7-1a06-40e6-bfe6-b511954e83f4.md"

    if ($xlfExt -eq "zh-cn") {
        $xlf45 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/45919da4390c41ba5f12730301a18a8385731388/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.zh-cn.xlf"
        $xlf30 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c61ea852603dcfc67c917c8b051a2d941426a724/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/304848e7-1a06-40e6-bfe6-b511954e83f4.bd312e9ba12e0136d876f519fa559c90eabf7655.zh-cn.xlf"
        $handoff45DateFile = "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.zh-cn.xlf"
        $handoff30DateFile = "304848e7-1a06-40e6-bfe6-b511954e83f4.bd312e9ba12e0136d876f519fa559c90eabf7655.zh-cn.xlf"
        $handoffDate45 = "2016-03-11 22:25:36"
        $handbackDate45 = "2016-03-11 22:25:50"
        $handoffDate30 = "2016-03-11 22:25:06"
    } else {
        $xlf45 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a67457509ea8308d9ccf5ffeb00ed10889f5f71f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.de-de.xlf"
        $xlf30 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f34ab406754316f5835a23911e15d26aab3807df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/304848e7-1a06-40e6-bfe6-b511954e83f4.bd312e9ba12e0136d876f519fa559c90eabf7655.de-de.xlf"
        $handoff45DateFile = "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.de-de.xlf"
        $handoff30DateFile = "304848e7-1a06-40e6-bfe6-b511954e83f4.bd312e9ba12e0136d876f519fa559c90eabf7655.de-de.xlf"
        $handoffDate45 = "2016-03-11 22:25:38"
        $handbackDate45 = "2016-03-11 22:25:56"
        $handoffDate30 = "2016-03-11 22:25:15"
    }

    # Clear existing hyperlinks for the two source-file rows before rebuilding.
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("B2").Hyperlinks.Delete()
    $ws.Range("D2").Hyperlinks.Delete()
    $ws.Range("A3").Hyperlinks.Delete()
    $ws.Range("B3").Hyperlinks.Delete()
    $ws.Range("D3").Hyperlinks.Delete()

    # --- Row 2: 45e61196 file, now handed back -----------------------------
    $ws.Range("A2").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.md"
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("D2").Value = $handoff45DateFile
    $ws.Range("E2").Value = $handoffDate45
    $ws.Range("F2").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.md"
    $ws.Range("G2").Value = $handoff45DateFile
    $ws.Range("H2").Value = $handbackDate45
    $ws.Range("I2").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl45, "", "", "45e61196-dd2b-4863-a39e-77c67d4820a6.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl45, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $xlf45, "", "", $handoff45DateFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl45, "", "", "45e61196-dd2b-4863-a39e-77c67d4820a6.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlf45, "", "", $handoff45DateFile) | Out-Null
    $ws.Range("F2").Font.Underline = 2
    $ws.Range("F2").Font.Color = 15570276
    $ws.Range("G2").Font.Underline = 2
    $ws.Range("G2").Font.Color = 15570276

    # --- Row 3: 304848e7 file, still ready for handoff ----------------------
    $ws.Range("A3").Value = "304848e7-1a06-40e6-bfe6-b511954e83f4.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $handoff30DateFile
    $ws.Range("E3").Value = $handoffDate30
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl30, "", "", "304848e7-1a06-40e6-bfe6-b511954e83f4.md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl30, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $xlf30, "", "", $handoff30DateFile) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
Fill-LangSheet $zh "zh-cn"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
Fill-LangSheet $de "de-de"
